# Update LR-pairs data with new TPM-derived values (Hbegf-Cd9)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 5.94498
$ws.Range("H2").Value = 17.83494
$ws.Range("I2").Value = 0.4679240463447598
$ws.Range("J2").Value = 0.4679240463447597
$ws.Range("M2").Value = 18.62071366666667
$ws.Range("N2").Value = 55.862141
$ws.Range("O2").Value = 0.1068221077965325
$ws.Range("P2").Value = 0.1068221077965325
$ws.Range("Q2").Value = 110.69977033406
$ws.Range("R2").Value = 996.29793300654
$ws.Range("S2").Value = 0.04998463291922958
$ws.Range("T2").Value = 0.04998463291922958

# Row 3
$ws.Range("G3").Value = 5.94498
$ws.Range("H3").Value = 17.83494
$ws.Range("I3").Value = 0.4679240463447598
$ws.Range("J3").Value = 0.4679240463447597
$ws.Range("O3").Value = 0.7040307798496723
$ws.Range("P3").Value = 0.7040307798496723
$ws.Range("Q3").Value = 729.5872291334599
$ws.Range("R3").Value = 6566.285062201139
$ws.Range("S3").Value = 0.3294329312585154
$ws.Range("T3").Value = 0.3294329312585154

# Row 4
$ws.Range("G4").Value = 5.94498
$ws.Range("H4").Value = 17.83494
$ws.Range("I4").Value = 0.4679240463447598
$ws.Range("J4").Value = 0.4679240463447597
$ws.Range("M4").Value = 32.97121066666667
$ws.Range("N4").Value = 98.91363200000001
$ws.Range("O4").Value = 0.1891471123537951
$ws.Range("P4").Value = 0.1891471123537951
$ws.Range("Q4").Value = 196.01318798912
$ws.Range("R4").Value = 1764.11869190208
$ws.Range("S4").Value = 0.08850648216701471
$ws.Range("T4").Value = 0.08850648216701468

# Row 5
$ws.Range("H5").Value = 8.352077
$ws.Range("I5").Value = 0.219128164447035
$ws.Range("J5").Value = 0.219128164447035
$ws.Range("M5").Value = 18.62071366666667
$ws.Range("N5").Value = 55.862141
$ws.Range("O5").Value = 0.1068221077965325
$ws.Range("P5").Value = 0.1068221077965325
$ws.Range("Q5").Value = 51.84054477965078
$ws.Range("R5").Value = 466.564903016857
$ws.Range("S5").Value = 0.02340773240381746
$ws.Range("T5").Value = 0.02340773240381746

# Row 6
$ws.Range("H6").Value = 8.352077
$ws.Range("I6").Value = 0.219128164447035
$ws.Range("J6").Value = 0.219128164447035
$ws.Range("O6").Value = 0.7040307798496723
$ws.Range("P6").Value = 0.7040307798496723
$ws.Range("R6").Value = 3074.981942381286
$ws.Range("S6").Value = 0.1542729725026732
$ws.Range("T6").Value = 0.1542729725026732

# Row 7
$ws.Range("H7").Value = 8.352077
$ws.Range("I7").Value = 0.219128164447035
$ws.Range("J7").Value = 0.219128164447035
$ws.Range("M7").Value = 32.97121066666667
$ws.Range("N7").Value = 98.91363200000001
$ws.Range("O7").Value = 0.1891471123537951
$ws.Range("P7").Value = 0.1891471123537951
$ws.Range("Q7").Value = 91.79269675707378
$ws.Range("R7").Value = 826.134270813664
$ws.Range("S7").Value = 0.04144745954054421
$ws.Range("T7").Value = 0.04144745954054421

# Row 8
$ws.Range("G8").Value = 3.976005
$ws.Range("H8").Value = 11.928015
$ws.Range("I8").Value = 0.3129477892082053
$ws.Range("J8").Value = 0.3129477892082053
$ws.Range("M8").Value = 18.62071366666667
$ws.Range("N8").Value = 55.862141
$ws.Range("O8").Value = 0.1068221077965325
$ws.Range("P8").Value = 0.1068221077965325
$ws.Range("Q8").Value = 74.03605064223501
$ws.Range("R8").Value = 666.324455780115
$ws.Range("S8").Value = 0.03342974247348543
$ws.Range("T8").Value = 0.03342974247348543

# Row 9
$ws.Range("G9").Value = 3.976005
$ws.Range("H9").Value = 11.928015
$ws.Range("I9").Value = 0.3129477892082053
$ws.Range("J9").Value = 0.3129477892082053
$ws.Range("O9").Value = 0.7040307798496723
$ws.Range("P9").Value = 0.7040307798496723
$ws.Range("Q9").Value = 487.948230434885
$ws.Range("R9").Value = 4391.534073913965
$ws.Range("S9").Value = 0.2203248760884837
$ws.Range("T9").Value = 0.2203248760884837

# Row 10
$ws.Range("G10").Value = 3.976005
$ws.Range("H10").Value = 11.928015
$ws.Range("I10").Value = 0.3129477892082053
$ws.Range("J10").Value = 0.3129477892082053
$ws.Range("M10").Value = 32.97121066666667
$ws.Range("N10").Value = 98.91363200000001
$ws.Range("O10").Value = 0.1891471123537951
$ws.Range("P10").Value = 0.1891471123537951
$ws.Range("Q10").Value = 131.09369846672
$ws.Range("R10").Value = 1179.84328620048
$ws.Range("S10").Value = 0.0591931706462362
$ws.Range("T10").Value = 0.05919317064623619
